# Adds three new grouping-variable rows (rice, extension_training,
# animal_ded_area) to the bottom of the table on Sheet1.
#
# Cell values are entered in the specific order below so that the
# workbook's shared-string table is built up in the same sequence as the
# authored edit (new strings are interned in first-use order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - rice
$ws.Range("A10").Value = "rice"
# Row 11 - extension_training
$ws.Range("A11").Value = "extension_training"
# Row 10 continued
$ws.Range("B10").Value = "Household cultivated rice"
$ws.Range("C10").Value = "Rice"
# Row 11 continued
$ws.Range("C11").Value = "Extension Training"
$ws.Range("B11").Value = "Household Received Extension Training"
# Row 12 - animal_ded_area
$ws.Range("B12").Value = "Household dedicated space to livestock"
$ws.Range("C12").Value = "Dedicated Livestock Area"
$ws.Range("A12").Value = "animal_ded_area"

# Shared Levels/Labels/level columns for the three new rows
$ws.Range("D10").Value = "0,1"
$ws.Range("E10").Value = "Yes,No"
$ws.Range("F10").Value = "All"

$ws.Range("D11").Value = "0,1"
$ws.Range("E11").Value = "Yes,No"
$ws.Range("F11").Value = "All"

$ws.Range("D12").Value = "0,1"
$ws.Range("E12").Value = "Yes,No"
$ws.Range("F12").Value = "All"

# Column C widened to fit the new, longer "Dedicated Livestock Area" label
$ws.Columns.Item(3).ColumnWidth = 27.59

# Active cell ends on the row below the newly-entered data, as it would
# after typing the last row and pressing Enter
[void]$ws.Range("A13").Select()
